$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures ---
# "Valor Mora" total (E11): 462800 -> 254800
$ws.Range("E11").Value2 = 254800
# "Cant. Trabajadores" (C13): 7 -> 4
$ws.Range("C13").Value2 = 4

# --- Copy the "last row" bottom-border formatting from the current last
#     data row (26) onto what will become the new last data row (22),
#     before we overwrite/remove any rows. ---
$ws.Range("B26:J26").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Rewrite the worker detail table (rows 16-22) with the updated data ---
# Row 16: LUIS GUILLERMO DELGADO VALENCIA - period 2407
$ws.Range("B16").Value2 = "CC"
$ws.Range("C16").Value2 = "1110572350"
$ws.Range("D16").Value2 = "LUIS GUILLERMO DELGADO VALENCIA"
$ws.Range("E16").Value2 = "2407"
$ws.Range("F16").Value2 = 52000
$ws.Range("G16").Value2 = 1300000

# Row 17: CARLOS AUGUSTO GONZALEZ GUTIERREZ - period 2505
$ws.Range("B17").Value2 = "CC"
$ws.Range("C17").Value2 = "1047377486"
$ws.Range("D17").Value2 = "CARLOS AUGUSTO GONZALEZ GUTIERREZ"
$ws.Range("E17").Value2 = "2505"
$ws.Range("F17").Value2 = 52000
$ws.Range("G17").Value2 = 1300000

# Row 18: DIEGO FELIPE GARCIA MONTEALEGRE - period 2505
$ws.Range("B18").Value2 = "CC"
$ws.Range("C18").Value2 = "14137062"
$ws.Range("D18").Value2 = "DIEGO FELIPE GARCIA MONTEALEGRE"
$ws.Range("E18").Value2 = "2505"
$ws.Range("F18").Value2 = 52000
$ws.Range("G18").Value2 = 1300000

# Row 19: JOAN SEBASTIAN VILLARREAL GARZON - period 2505
$ws.Range("B19").Value2 = "CC"
$ws.Range("C19").Value2 = "1001343476"
$ws.Range("D19").Value2 = "JOAN SEBASTIAN VILLARREAL GARZON"
$ws.Range("E19").Value2 = "2505"
$ws.Range("F19").Value2 = 52000
$ws.Range("G19").Value2 = 1300000

# Row 20: CARLOS AUGUSTO GONZALEZ GUTIERREZ - period 2506
$ws.Range("B20").Value2 = "CC"
$ws.Range("C20").Value2 = "1047377486"
$ws.Range("D20").Value2 = "CARLOS AUGUSTO GONZALEZ GUTIERREZ"
$ws.Range("E20").Value2 = "2506"
$ws.Range("F20").Value2 = 15600
$ws.Range("G20").Value2 = 1300000

# Row 21: DIEGO FELIPE GARCIA MONTEALEGRE - period 2506
$ws.Range("B21").Value2 = "CC"
$ws.Range("C21").Value2 = "14137062"
$ws.Range("D21").Value2 = "DIEGO FELIPE GARCIA MONTEALEGRE"
$ws.Range("E21").Value2 = "2506"
$ws.Range("F21").Value2 = 15600
$ws.Range("G21").Value2 = 1300000

# Row 22: JOAN SEBASTIAN VILLARREAL GARZON - period 2506 (last data row)
$ws.Range("B22").Value2 = "CC"
$ws.Range("C22").Value2 = "1001343476"
$ws.Range("D22").Value2 = "JOAN SEBASTIAN VILLARREAL GARZON"
$ws.Range("E22").Value2 = "2506"
$ws.Range("F22").Value2 = 15600
$ws.Range("G22").Value2 = 1300000

# --- Remove the now-unused trailing data rows (old rows 23-26), which
#     also shifts the signature block footer up from rows 31/32 to 27/28 ---
$ws.Rows("23:26").Delete()
